$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.06578863049612835
$ws.Range("C2").Value = 0.5873467334423524
$ws.Range("D2").Value = 0.6636964401451462
$ws.Range("E2").Value = 0.8146756656149404
$ws.Range("F2").Value = 0.8199372453822267

$ws.Range("B3").Value = -0.02523324468164586
$ws.Range("C3").Value = 0.505465297153278
$ws.Range("D3").Value = 0.4821656161413891
$ws.Range("E3").Value = 0.6943814629880244
$ws.Range("F3").Value = 0.7008277088516904

$ws.Range("B4").Value = 0.005608450684684903
$ws.Range("C4").Value = 0.4847052646643364
$ws.Range("D4").Value = 0.4358015571763177
$ws.Range("E4").Value = 0.6601526771712115
$ws.Range("F4").Value = 0.6668308405092866

$ws.Range("B5").Value = -0.01284610207697986
$ws.Range("C5").Value = 0.4921832099637495
$ws.Range("D5").Value = 0.500471272416235
$ws.Range("E5").Value = 0.7074399426214462
$ws.Range("F5").Value = 0.7161102612933321
$ws.Range("G5").Value = 41

$ws.Range("B6").Value = -0.09768167243710632
$ws.Range("C6").Value = 0.5114466272893242
$ws.Range("D6").Value = 0.4882265457959837
$ws.Range("E6").Value = 0.6987320987302528
$ws.Range("F6").Value = 0.7033071860062781
$ws.Range("G6").Value = 31

$ws.Range("B7").Value = -0.1103781127069636
$ws.Range("C7").Value = 0.5392641884755303
$ws.Range("D7").Value = 0.5358859154982463
$ws.Range("E7").Value = 0.7320422907853387
$ws.Range("F7").Value = 0.7364823498096339
$ws.Range("G7").Value = 29

$ws.Range("B8").Value = -0.09365692653440394
$ws.Range("C8").Value = 0.5420077535482939
$ws.Range("D8").Value = 0.5495523960454577
$ws.Range("E8").Value = 0.7413180127620383
$ws.Range("F8").Value = 0.7493864402156157
$ws.Range("G8").Value = 27

$ws.Range("B9").Value = -0.0615712106401871
$ws.Range("C9").Value = 0.5403710845059494
$ws.Range("D9").Value = 0.5803707086519223
$ws.Range("E9").Value = 0.7618206538627856
$ws.Range("F9").Value = 0.7801358214643087
$ws.Range("G9").Value = 19

$ws.Range("B10").Value = -0.01679119764512867
$ws.Range("C10").Value = 0.4453590725437407
$ws.Range("D10").Value = 0.4555235998480789
$ws.Range("E10").Value = 0.6749248845968556
$ws.Range("F10").Value = 0.704717858846984
$ws.Range("G10").Value = 12

$ws.Range("B11").Value = 0.3510802739590448
$ws.Range("C11").Value = 0.6628255016558736
$ws.Range("D11").Value = 0.9409865269452651
$ws.Range("E11").Value = 0.970044600492815
$ws.Range("F11").Value = 1.011020009805757
